# Estimation.xlsx — "added the binding object that was missing."
#
# The MySQL-for-Excel "Edit MySQL Data" binding had been dropped, leaving
# the estimate next to "Create Questions" (row 5) blank. Restore the
# missing value and the (hidden) helper defined name that the
# MySQL-for-Excel add-in's date/time formatting relies on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing estimate (row 5, col C).
# Selecting the cell afterwards mirrors where the user's cursor was left.
$ws.Range("C5").Value = "2 hrs"
$ws.Range("C5").Select() | Out-Null

# Re-add the hidden workbook-scoped name the MySQL-for-Excel binding
# object depends on for localized date/time formatting.
$formula = "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)"

$existing = $null
foreach ($n in $wb.Names) {
    if ($n.Name -eq "LOCAL_MYSQL_DATE_FORMAT") { $existing = $n }
}
if ($existing -ne $null) {
    $existing.Delete()
}

$name = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $formula)
$name.Visible = $false
